$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet: a new (blank) column was inserted before the
# old "Late" column (was column N, 14), shifting the old N/O/P
# ("Late" / heading / "Outstanding") columns one to the right (O/P/Q).
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns.Item(14).Insert() | Out-Null

# The newly inserted column kept a custom (non bestFit) width, matching the
# width used elsewhere in the sheet.
$wsSchedule.Columns.Item(14).ColumnWidth = 9.8

# Repayment schedule becomes the active sheet/tab, with K15 selected
# (this also clears the "Transactions" sheet's previous tab-selected flag).
$wsSchedule.Range("K15").Select() | Out-Null
$wsSchedule.Activate() | Out-Null
